$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2, shifting the existing data down.
$ws.Rows(2).Insert()

# Populate the new row with the new group/player correspondence.
$ws.Range("A2").Value = "tzorec"
$ws.Range("B2").Value = 120
$ws.Range("C2").Value = 4
$ws.Range("D2").Value = 2

# Match the author's final selection/view state.
$ws.Range("D2").Select()
